$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.886.79"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.452.50"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.446.01"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.51%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "4.044.88"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.96%  "
$ws.Range("D17").Value = "64.958.81"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "3.458.49"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").Value = "2.909.82"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0749"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.787"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0311"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "320.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "
